$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 changes from "001" to "002" (text, not a number) - force text then clear
# the number-format override so no extra style is introduced on the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").ClearFormats()

# N2: report date text changes
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric financial figures for row 2
$ws.Range("O2").Value = 562078648.1
$ws.Range("P2").Value = 87549711.44
$ws.Range("Q2").Value = 64227611.37
$ws.Range("S2").Value = 160109216.78
$ws.Range("U2").Value = 128079354.25
$ws.Range("W2").Value = 254762296.94
$ws.Range("X2").Value = 146302260.58
$ws.Range("Z2").Value = 206026.61
$ws.Range("AB2").Value = 307316351.16
$ws.Range("AF2").Value = 188.8014555444
$ws.Range("AG2").Value = 45.3250266313
